$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) cells whose new value parses as a plain number need the
# cell pre-formatted as Text, otherwise Excel silently coerces the
# assigned string into a numeric value (dropping trailing zeros, etc).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.66"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5133"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3980"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08455"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.83"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.283"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.289"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.015"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001112"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.53"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06768"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.79"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.970"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.271"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.08"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.76"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.393"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.33"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.801"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.634"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02436"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06507"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.966"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.267"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.194"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6450"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.052"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.23"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.012"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6084"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.14"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.007"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.207"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.68"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.210"

# Remaining Price (D) cells (non-numeric-looking strings) and all
# Volume(1h) (E) percentage-strings can be assigned directly as text.
$ws.Range("D2").Value = "28.563.83"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.886.78"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.99%  "
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "1.889.07"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("E21").Value = "  -1.34%  "
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "28.627.37"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "2.097.94"
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  -6.94%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  -1.95%  "
